# Applies the "Atualizacao de bases das ligas" results update for the
# "Germany Verbandsliga" sheet:
#   - 5 pairs of adjacent rows had their match data (everything but the
#     running id in column A) swapped between the two rows;
#   - 2 new match result rows are appended at the bottom of the table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$matchCols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB")

function Set-MatchRow {
    param($r, $vals)
    for ($i = 0; $i -lt $matchCols.Length; $i++) {
        $ws.Range($matchCols[$i] + $r).Value = $vals[$i]
    }
}

# --- 1) Fix the 5 row pairs whose home/away data was swapped ---
Set-MatchRow 16 @(14, 7138607, 'Germany Verbandsliga', 45168.625, 'Rot Weiss Walldorf II', 'Turnerschaft OberRoden', 3, 2, 'H', 2.25, 3.75, 2.5, 2.25, 3.8, 2.45, 0, 1.8, 2, 3.75, 1.95, 1.85, 1.25, -1, -1, 0.8, -1, 0.95, -1)
Set-MatchRow 17 @(15, 7138608, 'Germany Verbandsliga', 45168.625, 'SV UnterFlockenbach', 'SC Dortelweil', 1, 1, 'D', 1.083, 9, 16, 1.125, 7.5, 13, -2.5, 1.775, 1.925, 4.25, 1.975, 1.825, -1, 6.5, -1, -1, 0.925, -1, 0.825)
Set-MatchRow 69 @(67, 7423699, 'Germany Verbandsliga', 45233.66666666666, 'SG 2000 MulheimKarlich', 'Ahrweiler BC', 2, 2, 'D', 2.2, 5, 2.2, 2.2, 4.75, 2.2, 0, 1.9, 1.9, 4.25, 1.775, 2.025, -1, 3.75, -1, 0, 0, -0.5, 0.5125)
Set-MatchRow 70 @(68, 7423700, 'Germany Verbandsliga', 45233.66666666666, 'TuS Hornau', 'FC Burgsolms', 3, 0, 'H', 1.727, 4.5, 3.2, 1.727, 4.5, 3.2, -0.5, 1.775, 2.025, 3.5, 1.85, 1.95, 0.7270000000000001, -1, -1, 0.7749999999999999, -1, -1, 0.95)
Set-MatchRow 86 @(84, 7511958, 'Germany Verbandsliga', 45254.66666666666, 'SpVgg EGC Wirges', 'SG 2000 MulheimKarlich', 2, 1, 'H', 4.333, 4, 1.571, 4.2, 4, 1.571, 1, 1.875, 1.925, 3.75, 1.925, 1.875, 3.2, -1, -1, 0.875, -1, -1, 0.875)
Set-MatchRow 87 @(85, 7511976, 'Germany Verbandsliga', 45254.66666666666, 'DJK Bad Homburg', 'SG Bornheim 1945 GrunWeiss', 4, 0, 'H', 2, 3.75, 2.9, 1.8, 4, 3.3, -0.5, 1.85, 1.95, 3.5, 1.975, 1.825, 0.8, -1, -1, 0.8500000000000001, -1, 0.9750000000000001, -1)
Set-MatchRow 117 @(115, 8013719, 'Germany Verbandsliga', 45378.66666666666, 'SV Pars NeuIsenburg', 'SG Bornheim 1945 GrunWeiss', 2, 3, 'A', 1.5, 4.75, 4.2, 1.444, 5, 4.75, -1.25, 1.85, 1.95, 4, 1.875, 1.925, -1, -1, 3.75, -1, 0.95, 0.875, -1)
Set-MatchRow 118 @(116, 8014741, 'Germany Verbandsliga', 45378.66666666666, 'FC Concordia 03', 'FV Preussen Eberswalde', 1, 1, 'D', 2, 4, 2.75, 2, 4, 2.8, -0.25, 1.825, 1.975, 4.25, 2, 1.8, -1, 3, -1, -0.5, 0.4875, -1, 0.8)
Set-MatchRow 143 @(141, 8121117, 'Germany Verbandsliga', 45403.41666666666, 'RotWeiss Frankfurt', 'FCA 04 Darmstadt', 2, 1, 'H', 2.7, 3.75, 2.1, 2.75, 3.75, 2.1, 0.25, 1.875, 1.925, 3.5, 1.875, 1.925, 1.75, -1, -1, 0.875, -1, -1, 0.925)
Set-MatchRow 144 @(142, 8121110, 'Germany Verbandsliga', 45403.41666666666, 'FC Astoria Walldorf II', 'SV Spielberg', 1, 1, 'D', 2, 3.6, 3, 2, 3.6, 3, -0.25, 1.825, 1.975, 3.25, 1.95, 1.85, -1, 2.6, -1, -0.5, 0.4875, -1, 0.8500000000000001)

# --- 2) Append the 2 new match rows, copying the row-144 cell format
#        (bold/bordered/centered id column, yyyy-mm-dd date column) ---
$ws.Range("A144:D144").Copy() | Out-Null
$ws.Range("A145:D146").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

Set-MatchRow 145 @(143, 8131039, 'Germany Verbandsliga', 45405.60416666666, 'FSV Offenbach', 'TuS Russingen', 1, 2, 'A', 1.85, 3.8, 3.25, 1.85, 3.8, 3.2, -0.5, 1.9, 1.9, 3.25, 1.85, 1.95, -1, -1, 2.2, -1, 0.8999999999999999, -0.5, 0.475)
Set-MatchRow 146 @(144, 8136061, 'Germany Verbandsliga', 45406.55208333334, '1 FC Lok Stendal', 'SSV 80 Gardelegen', 3, 1, 'H', 2.4, 3.5, 2.4, 2.4, 3.5, 2.4, 0, 1.9, 1.9, 3, 1.95, 1.85, 1.4, -1, -1, 0.8999999999999999, -1, 0.95, -1)

$ws.Range("A1:AB146").Select() | Out-Null
